$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C gets as wide as column B ---
$ws.Columns.Item(3).ColumnWidth = $ws.Columns.Item(2).ColumnWidth

# --- Header cell C2 gets a new "light" variant of the blue header style: ---
# same bold/size/color font as the other header cells, but a light-blue fill.
$ws.Range("B2").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("C2").Font.Size = 16
$ws.Range("C2").Font.Bold = $true
$ws.Range("C2").Font.Color = 16777215
$ws.Range("C2").Interior.Color = 16642787

# --- Existing row 3 (A3) changes its value, new cells are added alongside it,
#     and a brand-new row 4 appears with another entry. All of these share
#     A3's original (unstyled) look. ---
$ws.Range("A3").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("A3").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)

$ws.Range("B3").Value = "Activo: Cuenta"
$ws.Range("C3").Value = "No existe un monto FOSEDE para el tipo de moneda"

# Numeric-looking labels must stay text, so build them as a text formula and
# then collapse the formula down to its literal value (keeps the existing
# cell format untouched, unlike typing the digits straight in which Excel
# would interpret as a number).
$ws.Range("A3").Formula = "=""18"""
$ws.Range("A3").Copy()
$ws.Range("A3").PasteSpecial(-4163)

$ws.Range("A4").Formula = "=""30"""
$ws.Range("A4").Copy()
$ws.Range("A4").PasteSpecial(-4163)
